$wb = $excel.ActiveWorkbook

# --- Add the new "Test results" worksheet, positioned right after the
# first sheet ("Accounts receivable ledger") so it becomes the 2nd tab. ---
$wsLedger = $wb.Worksheets.Item(1)

$newSheet = $wb.Worksheets.Add($null, $wsLedger)
$newSheet.Name = "Test results"

# Re-resolve these by name AFTER the insertion shifts sheet positions.
$wsCustomers = $wb.Worksheets.Item("Customers")
$wsPayments = $wb.Worksheets.Item("Payments")

# Column widths (A:B wide, C medium) matching the authored layout.
$newSheet.Columns.Item(1).ColumnWidth = 19.6667
$newSheet.Columns.Item(2).ColumnWidth = 19.6667
$newSheet.Columns.Item(3).ColumnWidth = 15.4167

# Number formats for the date column (B) and accuracy/percentage column (E)
# set up-front so the written values pick up the right style without the
# engine inventing an ad-hoc custom format.
$newSheet.Range("B3:B7").NumberFormat = "d-mmm-yy"
$newSheet.Range("E3:E6").NumberFormat = "0%"

# Header row
$newSheet.Range("A1").Value = "Workflow type"
$newSheet.Range("C1").Value = "Model"
$newSheet.Range("D1").Value = "Time"
$newSheet.Range("B1").Value = "Instance"

# Totals / summary row
$newSheet.Range("D2").Value = "228m"

# Type 5 test row
$newSheet.Range("A3").Value = "Type 5"
$newSheet.Range("B3").Value = 46027
$newSheet.Range("D3").Value = "227m"
$newSheet.Range("E1").Value = "Accuracy"
$newSheet.Range("C3").Value = "Qwen3:8b"
$newSheet.Range("E3").Value = 1

# Type 1 test rows
$newSheet.Range("A4").Value = "Type 1"
$newSheet.Range("B4").Value = 46027
$newSheet.Range("C4").Value = "llama3.1:8bn"
$newSheet.Range("D4").Value = "5m"
$newSheet.Range("E4").Value = 0

$newSheet.Range("A5").Value = "Type 1"
$newSheet.Range("B5").Value = 46027
$newSheet.Range("C5").Value = "qwen3:8bn"
$newSheet.Range("D5").Value = "83m"
$newSheet.Range("E5").Value = 0

$newSheet.Range("A6").Value = "Type 1"
$newSheet.Range("B6").Value = 46027
$newSheet.Range("C6").Value = "deepseek-r1:14b"
$newSheet.Range("D6").Value = "32m"
$newSheet.Range("E6").Value = 0

$newSheet.Range("A7").Value = "Type 1"
$newSheet.Range("B7").Value = 46027
$newSheet.Range("C7").Value = "gpt-oss:20b"

# --- Update the remembered selections on each sheet. The last sheet
# selected becomes the active tab, so "Test results" is selected last. ---
$wsLedger.Range("D13").Select()
$wsCustomers.Range("G2").Select()
$wsPayments.Range("D6").Select()
$newSheet.Range("D7").Select()
